$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# --- Paragraph 1: "data Tag : CLabel -> CEnum -> Type where  <TAB>TZ : Tag l (l :: e)" ---

# "CLabel" (chars 12-17) and the trailing space (char 18):
# drop the lumMod/lumOff shading, keep plain accent2 scheme color.
$tr.Characters(12, 7).Font.Color.ObjectThemeColor = 6   # msoThemeColorAccent2

# "CEnum" (chars 22-26) and the trailing space (char 27):
# recolor from the D99694 tint to the solid C0504D red.
$tr.Characters(22, 6).Font.Color.RGB = 5066944           # 0xC0504D -> R+G*256+B*65536

# " :: " (chars 57-60) becomes " " (unstyled) + ":: " (accent3 scheme color)
$tr.Characters(58, 3).Font.Color.ObjectThemeColor = 7    # msoThemeColorAccent3

# --- Paragraph 2: "<TAB>TS : Tag l e -> Tag l (l' :: e)" ---

# ":: " (chars 91-93) becomes "::" (9BBB59 green) + " " (unstyled)
$tr.Characters(91, 2).Font.Color.RGB = 5880731            # 0x9BBB59 -> R+G*256+B*65536
